# Adding new progress as of date 04-Nov-2025.
#
# For every data row (3-26) on the "Training Dashboard" sheet:
#   - column H ("PERIOD TO EXPIRE") decreases by one day
#   - column I ("LAST UPDATE") is refreshed from "03-Nov-2025" to "04-Nov-2025"
#
# Column I holds the date as literal text (not a real Excel date), so the
# assignment is done with a leading apostrophe to stop Excel's automatic
# date recognition from turning the literal into a date serial number. The
# apostrophe/quote-prefix leaves a "number stored as text" style behind, so
# the original cell formatting is restored right after by copying the
# format from the neighbouring (unaffected) column J cell on the same row,
# which carries the identical base style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 26; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE
    $iCell.Value = "'" + "04-Nov-2025"

    # Re-apply the row's normal style so the quote-prefix used above to
    # keep the date as text doesn't leave a different cell style behind.
    $fmtSource = $ws.Cells.Item($row, 10)   # column J - same style, untouched
    $fmtSource.Copy()
    $iCell.PasteSpecial(-4122)              # xlPasteFormats
}

$excel.CutCopyMode = $false
